# Collapse the "от <tab> г." date placeholder into a single run whose
# text is a fill-in-the-blank line: "от____________________________20____г."
#
# In the original document this line is built from four runs:
#   1) "от "            (sz 22)
#   2) <tab>             (sz 20, en-GB)
#   3) " "               (sz 22)
#   4) "г."              (sz 22)
# The edit merges all of that into one run, keeping the first run's
# formatting (sz 22 / szCs 22), with the literal text
# "от____________________________20____г."

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "от ^t г.",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "от____________________________20____г.",
    2
) | Out-Null
